$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: the title field "${local_directorate}" becomes
# "${local_directorate_genitive}" -- two new runs ("_" and "genitive") are
# inserted right after the existing "directorate" run (and before the
# closing "}" run) inside the very first occurrence of this placeholder
# (the document "Θέμα:" / title paragraph).
# ---------------------------------------------------------------------------

# Locate the (first) "directorate" run that sits right after "local_" and
# right before the closing "}" of the ${local_directorate} placeholder.
$dirRng = $d.Content
$dirRng.Find.Execute("directorate", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dirStart = $dirRng.Start
$dirEnd = $dirRng.End

# Grab a formatted copy of the "_" run (the underscore between "local" and
# "directorate") and of the "directorate" run itself, so the new runs can
# reuse the exact same run formatting (rFonts/sz/szCs[/lang]).
$underscoreSrc = $d.Range($dirStart - 1, $dirStart)
$ftUnderscore = $underscoreSrc.FormattedText

$directorateSrc = $d.Range($dirStart, $dirEnd)
$ftDirectorate = $directorateSrc.FormattedText

# Insert a formatted copy of "_" right after "directorate".
$insUnderscore = $d.Range($dirEnd, $dirEnd)
$insUnderscore.FormattedText = $ftUnderscore

# Re-stamp its text (round-trip through a different value) so the new run
# is emitted as a clean, standalone <w:r> (no stale rsid carried over from
# the source run it was cloned from).
$underscoreNew = $d.Range($dirEnd, $dirEnd + 1)
$underscoreNew.Text = "#"
$underscoreNew = $d.Range($dirEnd, $dirEnd + 1)
$underscoreNew.Text = "_"

# Insert a formatted copy of "directorate" right after the new "_" run,
# then turn its text into "genitive" (same trick: this yields a clean new
# <w:r> that keeps the "directorate" run's formatting, including lang).
$genitiveStart = $dirEnd + 1
$insGenitive = $d.Range($genitiveStart, $genitiveStart)
$insGenitive.FormattedText = $ftDirectorate

$genitiveNew = $d.Range($genitiveStart, $genitiveStart + 11)
$genitiveNew.Text = "genitive"

# ---------------------------------------------------------------------------
# Change 2: drop the stray space before the comma in
# " για τους κάτωθι εκπαιδευτικούς , ως εξής:" ->
# " για τους κάτωθι εκπαιδευτικούς, ως εξής:"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("κάτωθι εκπαιδευτικούς , ως εξής:", $true, $false, $false, $false, $false, $true, 1, $false, "κάτωθι εκπαιδευτικούς, ως εξής:", 2) | Out-Null
